# Script results were refreshed: 2025-04-16 no longer has a page number
# ("NA"), and a new day (2025-04-17) was appended with no matching results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C51 ("NA") is cleared, matching the other "nothing found" rows above it.
$ws.Range("C51").ClearContents()

# Append the new row for 2025-04-17.
# The leading apostrophe forces the date-like string into literal text
# (same as the other Date-column cells) instead of Excel auto-converting
# it to a date serial number.
$ws.Range("A52").Value = "'2025-04-17"
$ws.Range("B52").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C52").Value = "NA"
$ws.Range("D52").Value = 1
